# Fix: scraped numbers (and a few text fields that happened to contain
# commas/periods) were mangled by the scraper using a naive
# "swap decimal comma<->dot" pass that didn't remove the thousands
# separator first. Re-apply the *correct* transform:
#   remove every "." (thousands separator) then turn "," into "."
# to every value that was bitten by it.
function Convert-ScrapedText([string]$s) {
    return $s.Replace(".", "").Replace(",", ".")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Razon social" / "Nombre Fantasia" entries mangled the same way ---
# Each group below is every cell sharing one bad value; update them together.
$textGroups = @(
    @("E63", "F63"),
    @("E69"),
    @("E70", "F70"),
    @("E72", "E125")
)

foreach ($group in $textGroups) {
    $old = $ws.Range($group[0]).Value()
    $new = Convert-ScrapedText $old
    foreach ($addr in $group) {
        $ws.Range($addr).Value = $new
    }
}

# --- "Importe" column (H2:H163): comma-decimal amounts like "95.169,60" ---
# These must stay text (so trailing zeros like ",60"/",00" survive), so force
# the cell to Text format before assigning, otherwise Excel/COM will parse the
# now dot-decimal-looking string as a number and round-trip it with floating
# point noise. Reset the cell style back to Normal afterwards so we don't
# leave a stray "Text" look to the cell beyond its format code.
for ($r = 2; $r -le 163; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value()
    $new = Convert-ScrapedText $old
    $cell.NumberFormat = "@"
    $cell.Value = $new
    $cell.Style = "Normal"
}
